$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 3.55
$ws.Range("I3").Value = 6.7
$ws.Range("K3").Value = 5.9
$ws.Range("M3").Value = 2.52
$ws.Range("N3").Value = 2.35
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.4
$ws.Range("U3").Value = 5.7
$ws.Range("W3").Value = 10
$ws.Range("Z3").Value = 5.9
$ws.Range("AA3").Value = 7.3
$ws.Range("AC3").Value = 175
$ws.Range("AD3").Value = 12.5
$ws.Range("AE3").Value = 40
$ws.Range("AF3").Value = 23
$ws.Range("AG3").Value = 175

# Row 4
$ws.Range("G4").Value = 2.27
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 3.35
$ws.Range("J4").Value = 1.11
$ws.Range("K4").Value = 5.5
$ws.Range("L4").Value = 1.47
$ws.Range("M4").Value = 2.5
$ws.Range("N4").Value = 2.4
$ws.Range("O4").Value = 1.5
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 2.02
$ws.Range("S4").Value = 1.72
$ws.Range("T4").Value = 6.1
$ws.Range("U4").Value = 10
$ws.Range("V4").Value = 9.25
$ws.Range("W4").Value = 23
$ws.Range("X4").Value = 21
$ws.Range("Z4").Value = 5.5
$ws.Range("AA4").Value = 5.7
$ws.Range("AB4").Value = 16.5
$ws.Range("AC4").Value = 100
$ws.Range("AD4").Value = 7.6
$ws.Range("AE4").Value = 16
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 37
$ws.Range("AI4").Value = 50

# Row 6
$ws.Range("G6").Value = 1.62
$ws.Range("H6").Value = 3.55
$ws.Range("I6").Value = 5.3
$ws.Range("L6").Value = 1.4
$ws.Range("M6").Value = 2.52
$ws.Range("N6").Value = 2.15
$ws.Range("O6").Value = 1.55
$ws.Range("P6").Value = 1.47
$ws.Range("Q6").Value = 2.32
$ws.Range("R6").Value = 2.15
$ws.Range("S6").Value = 1.55
$ws.Range("T6").Value = 5.2
$ws.Range("U6").Value = 6.3
$ws.Range("V6").Value = 8.75
$ws.Range("W6").Value = 11.25
$ws.Range("X6").Value = 15.5
$ws.Range("Y6").Value = 40
$ws.Range("Z6").Value = 7.6
$ws.Range("AA6").Value = 7.1
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 150
$ws.Range("AD6").Value = 11
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 19
$ws.Range("AG6").Value = 110
$ws.Range("AH6").Value = 70
$ws.Range("AI6").Value = 90

# Row 7
$ws.Range("G7").Value = 1.88
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 3.85
$ws.Range("L7").Value = 1.33
$ws.Range("M7").Value = 2.8
$ws.Range("N7").Value = 1.98
$ws.Range("O7").Value = 1.65
$ws.Range("P7").Value = 1.44
$ws.Range("Q7").Value = 2.42
$ws.Range("R7").Value = 1.82
$ws.Range("S7").Value = 1.78
$ws.Range("T7").Value = 6.5
$ws.Range("U7").Value = 8.5
$ws.Range("V7").Value = 8.5
$ws.Range("W7").Value = 16
$ws.Range("X7").Value = 16
$ws.Range("Y7").Value = 30
$ws.Range("Z7").Value = 9
$ws.Range("AA7").Value = 6.5
$ws.Range("AB7").Value = 16
$ws.Range("AC7").Value = 80
$ws.Range("AD7").Value = 10.25
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 13
$ws.Range("AG7").Value = 60
$ws.Range("AH7").Value = 37
$ws.Range("AI7").Value = 45
$ws.Range("AJ7").Value = 700

# Row 13
$ws.Range("G13").Value = 3.4
$ws.Range("I13").Value = 2.3
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.5
$ws.Range("T13").Value = 8.5
$ws.Range("U13").Value = 15
$ws.Range("V13").Value = 13
$ws.Range("AE13").Value = 10
$ws.Range("AG13").Value = 21
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 34
$ws.Range("AJ13").Value = 1250

# Row 16
$ws.Range("H16").Value = 3.4
$ws.Range("I16").Value = 3.8
$ws.Range("P16").Value = 1.5
$ws.Range("Q16").Value = 2.5
$ws.Range("T16").Value = 5.5
$ws.Range("U16").Value = 8.5
$ws.Range("V16").Value = 8.5
$ws.Range("W16").Value = 17
$ws.Range("X16").Value = 19
$ws.Range("Y16").Value = 29
$ws.Range("AA16").Value = 7
$ws.Range("AH16").Value = 34

# Row 17
$ws.Range("G17").Value = 2.8
$ws.Range("I17").Value = 2.5
$ws.Range("T17").Value = 7.5
$ws.Range("U17").Value = 13
$ws.Range("W17").Value = 29
$ws.Range("X17").Value = 26
$ws.Range("AC17").Value = 51
$ws.Range("AE17").Value = 11
$ws.Range("AF17").Value = 10
$ws.Range("AG17").Value = 23
$ws.Range("AI17").Value = 34
